$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update selections on the existing 3 sheets
# ---------------------------------------------------------------------
$wsMst    = $wb.Worksheets.Item("RAIIN_KBN_MST")
$wsDetail = $wb.Worksheets.Item("RAIIN_KBN_DETAIL")
$wsInf    = $wb.Worksheets.Item("RAIIN_KBN_INF")

$wsMst.Range("J23").Select() | Out-Null
$wsDetail.Range("C2").Select() | Out-Null
$wsInf.Range("E2").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. Add the 3 new sheets at the end of the workbook, in order
# ---------------------------------------------------------------------
$wsKoui     = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsKoui.Name = "RAIIN_KBN_KOUI"

$wsKouiMst  = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsKouiMst.Name = "KOUI_KBN_MST"

$wsItem     = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsItem.Name = "RAIIN_KBN_ITEM"

# ---------------------------------------------------------------------
# 3. RAIIN_KBN_KOUI (sheet4)
# ---------------------------------------------------------------------
$wsKoui.Range("A1").Value = "HP_ID"
$wsKoui.Range("B1").Value = "GRP_ID"
$wsKoui.Range("C1").Value = "KBN_CD"
$wsKoui.Range("D1").Value = "SEQ_NO"
$wsKoui.Range("E1").Value = "KOUI_KBN_ID"

$wsKoui.Range("A2").Value = 1
$wsKoui.Range("B2").Value = 9999
$wsKoui.Range("C2").Value = 999
$wsKoui.Range("D2").Value = 0
$wsKoui.Range("E2").Value = 999

$wsKoui.Range("G2").NumberFormat = "mm:ss.0"
$wsKoui.Range("J2").NumberFormat = "mm:ss.0"

$wsKoui.Columns.Item(5).ColumnWidth = 13.140625

$wsKoui.Range("C2").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. KOUI_KBN_MST (sheet5)
# ---------------------------------------------------------------------
$wsKouiMst.Range("A1").Value = "HP_ID"
$wsKouiMst.Range("B1").Value = "KOUI_KBN_ID"
$wsKouiMst.Range("C1").Value = "SORT_NO"
$wsKouiMst.Range("D1").Value = "KOUI_KBN1"
$wsKouiMst.Range("E1").Value = "KOUI_KBN2"
$wsKouiMst.Range("F1").Value = "KOUI_GRP_NAME"
$wsKouiMst.Range("G1").Value = "KOUI_NAME"

$wsKouiMst.Range("A2").Value = 1
$wsKouiMst.Range("B2").Value = 999
$wsKouiMst.Range("C2").Value = 1
$wsKouiMst.Range("D2").Value = 13
$wsKouiMst.Range("E2").Value = 13

# G2 must be written before F2 so the shared-string table gets
# KOUI_NAME_CHECK (23) before KOUI_GRP_NAME_CHECK (24)
$wsKouiMst.Range("G2").Value = "KOUI_NAME_CHECK"
$wsKouiMst.Range("F2").Value = "KOUI_GRP_NAME_CHECK"

$wsKouiMst.Columns.Item(1).ColumnWidth = 6.28515625
$wsKouiMst.Columns.Item(2).ColumnWidth = 13.140625
$wsKouiMst.Columns.Item(4).ColumnWidth = 11.140625
$wsKouiMst.Columns.Item(5).ColumnWidth = 11.140625
$wsKouiMst.Columns.Item(6).ColumnWidth = 23.140625
$wsKouiMst.Columns.Item(7).ColumnWidth = 18.85546875

$wsKouiMst.Range("H9").Select() | Out-Null

# ---------------------------------------------------------------------
# 5. RAIIN_KBN_ITEM (sheet6)
# ---------------------------------------------------------------------
$wsItem.Range("A1").Value = "HP_ID"
$wsItem.Range("B1").Value = "GRP_ID"
$wsItem.Range("C1").Value = "KBN_CD"
$wsItem.Range("D1").Value = "SEQ_NO"
$wsItem.Range("E1").Value = "ITEM_CD"
$wsItem.Range("F1").Value = "IS_EXCLUDE"
$wsItem.Range("G1").Value = "SORT_NO"

$wsItem.Range("A2").Value = 1
$wsItem.Range("B2").Value = 9999
$wsItem.Range("C2").Value = 999
$wsItem.Range("D2").Value = 0
$wsItem.Range("E2").Value = 613120001
$wsItem.Range("F2").Value = 0
$wsItem.Range("G2").Value = 0

$wsItem.Columns.Item(3).ColumnWidth = 8.140625
$wsItem.Columns.Item(4).ColumnWidth = 8.28515625
$wsItem.Columns.Item(5).ColumnWidth = 10
$wsItem.Columns.Item(6).ColumnWidth = 11.28515625

# RAIIN_KBN_ITEM ends up the active sheet (activeTab=5, tabSelected=1)
$wsItem.Range("E5").Select() | Out-Null
